$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Part Page" -> "Link"
$ws.Range("G1").Value = "Link"

# Connectors J11/J12, J2/J3, J9 quantity changed from "dni" (text) to numeric 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 0

# Remove stray note "stocked at mousser" from X2 row
$ws.Range("G28").Value = ""

# Restore a plain single-cell selection (was A28:XFD29)
$ws.Range("B1").Select()
